$d = $word.ActiveDocument

# 1) "have" -> "having" (highpoints geographical makeup range sentence)
$null = $d.Content.Find.Execute(
    "highpoints have such a wide range",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "highpoints having such a wide range", 2)

# 2) "Similarly " -> "Also " (start of paragraph)
$null = $d.Content.Find.Execute(
    "Similarly one might also incorrectly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Also one might also incorrectly", 2)

# 3) "There again" -> "Then again" (start of paragraph)
$null = $d.Content.Find.Execute(
    "There again what exactly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Then again what exactly", 2)

# 4) Insert "the summit of " before "Rainier is close" and move the
#    "_GoBack" bookmark so it sits right before "Rainier" (it currently
#    sits right after the closing "]" near the end of the equation line).
$null = $d.Content.Find.Execute(
    "average summer month temperature on Rainier is close",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "average summer month temperature on the summit of Rainier is close", 2)

$text = $d.Content.Text
$idx = $text.IndexOf("Rainier is close to or below the 32 F freezing point")
$bookmarkRange = $d.Range($idx, $idx)
$null = $d.Bookmarks.Add("_GoBack", $bookmarkRange)
